# Weekly update: insert the newest week's Ajo (garlic) price record at the
# top of the historical table (row 177), pushing all existing rows down by
# one. The oldest historical row therefore moves from row 208 to row 209.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 177; Excel shifts rows 177..208 down to
# 178..209 and carries formatting down from the row above (so column D keeps
# its date style).
$ws.Rows.Item(177).Insert()

# Populate the new row 177 with the latest week's record.
$ws.Cells.Item(177, 1).Value = 7
$ws.Cells.Item(177, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(177, 3).Value = "Ñuble"
$ws.Cells.Item(177, 4).Value = 44637
$ws.Cells.Item(177, 5).Value = 16
$ws.Cells.Item(177, 6).Value = 100112003
$ws.Cells.Item(177, 7).Value = "Ajo"
$ws.Cells.Item(177, 8).Value = "Chino"
$ws.Cells.Item(177, 9).Value = "Primera"
$ws.Cells.Item(177, 10).Value = 60
$ws.Cells.Item(177, 11).Value = 19000
$ws.Cells.Item(177, 12).Value = 20000
$ws.Cells.Item(177, 13).Value = 19500
$ws.Cells.Item(177, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(177, 15).Value = "China"
$ws.Cells.Item(177, 16).Value = 1950
$ws.Cells.Item(177, 17).Value = 10
$ws.Cells.Item(177, 18).Value = "Hortaliza"
